# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E16:E21) on Hoja1 is re-stated in ascending
# chronological order (2017-09 .. 2018-02) instead of the previous
# descending order, as part of refreshing the account-statement database.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "1709"
$ws.Range("E17").Value = "1710"
$ws.Range("E18").Value = "1711"
$ws.Range("E19").Value = "1712"
$ws.Range("E20").Value = "1801"
$ws.Range("E21").Value = "1802"
